$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 19.872027
$ws.Range("H2").Value = 59.61608099999999
$ws.Range("I2").Value = 0.2756064822985579
$ws.Range("J2").Value = 0.2756064822985579
$ws.Range("M2").Value = 211.980367
$ws.Range("N2").Value = 635.9411009999999
$ws.Range("O2").Value = 0.9885149156420702
$ws.Range("P2").Value = 0.9885149156420702
$ws.Range("Q2").Value = 4212.479576493908
$ws.Range("R2").Value = 37912.31618844517
$ws.Range("S2").Value = 0.2724411185997667
$ws.Range("T2").Value = 0.2724411185997667

$ws.Range("G3").Value = 19.872027
$ws.Range("H3").Value = 59.61608099999999
$ws.Range("I3").Value = 0.2756064822985579
$ws.Range("J3").Value = 0.2756064822985579
$ws.Range("O3").Value = 0.003992992409159323
$ws.Range("P3").Value = 0.003992992409159324
$ws.Range("Q3").Value = 17.015827183299
$ws.Range("R3").Value = 153.142444649691
$ws.Range("S3").Value = 0.001100494591733245
$ws.Range("T3").Value = 0.001100494591733245

$ws.Range("G4").Value = 19.872027
$ws.Range("H4").Value = 59.61608099999999
$ws.Range("I4").Value = 0.2756064822985579
$ws.Range("J4").Value = 0.2756064822985579
$ws.Range("O4").Value = 0.007492091948770576
$ws.Range("P4").Value = 0.007492091948770576
$ws.Range("Q4").Value = 31.926968242974
$ws.Range("R4").Value = 287.342714186766
$ws.Range("S4").Value = 0.002064869107058006
$ws.Range("T4").Value = 0.002064869107058006

$ws.Range("I5").Value = 0.1760995803479087
$ws.Range("J5").Value = 0.1760995803479087
$ws.Range("M5").Value = 211.980367
$ws.Range("N5").Value = 635.9411009999999
$ws.Range("O5").Value = 0.9885149156420702
$ws.Range("P5").Value = 0.9885149156420702
$ws.Range("Q5").Value = 2691.576335425675
$ws.Range("R5").Value = 24224.18701883107
$ws.Range("S5").Value = 0.1740770618122169
$ws.Range("T5").Value = 0.174077061812217

$ws.Range("I6").Value = 0.1760995803479087
$ws.Range("J6").Value = 0.1760995803479087
$ws.Range("O6").Value = 0.003992992409159323
$ws.Range("P6").Value = 0.003992992409159324
$ws.Range("S6").Value = 0.0007031642875853417
$ws.Range("T6").Value = 0.000703164287585342

$ws.Range("I7").Value = 0.1760995803479087
$ws.Range("J7").Value = 0.1760995803479087
$ws.Range("O7").Value = 0.007492091948770576
$ws.Range("P7").Value = 0.007492091948770576
$ws.Range("S7").Value = 0.001319354248106444
$ws.Range("T7").Value = 0.001319354248106444

$ws.Range("I8").Value = 0.5482939373535334
$ws.Range("J8").Value = 0.5482939373535334
$ws.Range("M8").Value = 211.980367
$ws.Range("N8").Value = 635.9411009999999
$ws.Range("O8").Value = 0.9885149156420702
$ws.Range("P8").Value = 0.9885149156420702
$ws.Range("Q8").Value = 8380.343574485205
$ws.Range("R8").Value = 75423.09217036683
$ws.Range("S8").Value = 0.5419967352300865
$ws.Range("T8").Value = 0.5419967352300865

$ws.Range("I9").Value = 0.5482939373535334
$ws.Range("J9").Value = 0.5482939373535334
$ws.Range("O9").Value = 0.003992992409159323
$ws.Range("P9").Value = 0.003992992409159324
$ws.Range("S9").Value = 0.002189333529840736
$ws.Range("T9").Value = 0.002189333529840736

$ws.Range("I10").Value = 0.5482939373535334
$ws.Range("J10").Value = 0.5482939373535334
$ws.Range("O10").Value = 0.007492091948770576
$ws.Range("P10").Value = 0.007492091948770576
$ws.Range("R10").Value = 571.642099334386
$ws.Range("S10").Value = 0.004107868593606126
$ws.Range("T10").Value = 0.004107868593606126
